$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Types")

$ws.Range("F29").Value = "BIT STRING (SIZE (20)"
$ws.Range("F30").Value = "BIT STRING (SIZE (28)"
$ws.Range("F31").Value = "BIT STRING (SIZE(18)"
$ws.Range("F32").Value = "BIT STRING (SIZE(21)"
$ws.Range("F33").Value = "BIT STRING (SIZE(20)"
$ws.Range("F34").Value = "BIT STRING (SIZE(18)"
$ws.Range("F35").Value = "BIT STRING (SIZE(21)"
$ws.Range("F36").Value = "BIT STRING (SIZE (22..32)"
$ws.Range("F37").Value = "BIT STRING (SIZE(22..32)"
